$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Cases" query (row 2 / cell B2) had an erroneous extra `Cohort` column
# appended to it (the `co:cohort` match was never meant to be selected here -
# that logic belongs to a different query). Remove the trailing
# `coalesce(co.cohort_description, '') AS `Cohort`` line so the query text
# matches the corrected/fixed query used elsewhere in the workbook.
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Flat-Coated Retriever'] 
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value2 = $newCasesQuery

# Move/refresh the active selection to B2 (was sitting on B4 previously).
$ws.Range("B2").Select() | Out-Null
